$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.5979736666666666
$ws.Range("H2").Value = 1.793921
$ws.Range("I2").Value = 0.03342655292740804
$ws.Range("J2").Value = 0.03342655292740804
$ws.Range("M2").Value = 28.19948866666667
$ws.Range("N2").Value = 84.598466
$ws.Range("O2").Value = 0.7357427920402423
$ws.Range("P2").Value = 0.7357427920402422
$ws.Range("Q2").Value = 16.86255163613178
$ws.Range("R2").Value = 151.762964725186
$ws.Range("S2").Value = 0.02459334537909213
$ws.Range("T2").Value = 0.02459334537909212

$ws.Range("G3").Value = 0.5979736666666666
$ws.Range("H3").Value = 1.793921
$ws.Range("I3").Value = 0.03342655292740804
$ws.Range("J3").Value = 0.03342655292740804
$ws.Range("O3").Value = 0.2029336910395279
$ws.Range("P3").Value = 0.2029336910395278
$ws.Range("Q3").Value = 4.651054527324111
$ws.Range("R3").Value = 41.859490745917
$ws.Range("S3").Value = 0.006783373764287049
$ws.Range("T3").Value = 0.006783373764287048

$ws.Range("G4").Value = 0.5979736666666666
$ws.Range("H4").Value = 1.793921
$ws.Range("I4").Value = 0.03342655292740804
$ws.Range("J4").Value = 0.03342655292740804
$ws.Range("M4").Value = 2.350402666666667
$ws.Range("N4").Value = 7.051208000000001
$ws.Range("O4").Value = 0.0613235169202299
$ws.Range("P4").Value = 0.06132351692022989
$ws.Range("Q4").Value = 1.405478900729778
$ws.Range("R4").Value = 12.649310106568
$ws.Range("S4").Value = 0.002049833784028867
$ws.Range("T4").Value = 0.002049833784028867

$ws.Range("I5").Value = 0.8874158839838097
$ws.Range("J5").Value = 0.8874158839838097
$ws.Range("M5").Value = 28.19948866666667
$ws.Range("N5").Value = 84.598466
$ws.Range("O5").Value = 0.7357427920402423
$ws.Range("P5").Value = 0.7357427920402422
$ws.Range("Q5").Value = 447.6709339098718
$ws.Range("R5").Value = 4029.038405188846
$ws.Range("S5").Value = 0.6529098401831078
$ws.Range("T5").Value = 0.6529098401831077

$ws.Range("I6").Value = 0.8874158839838097
$ws.Range("J6").Value = 0.8874158839838097
$ws.Range("O6").Value = 0.2029336910395279
$ws.Range("P6").Value = 0.2029336910395278
$ws.Range("S6").Value = 0.1800865808239399
$ws.Range("T6").Value = 0.1800865808239399

$ws.Range("I7").Value = 0.8874158839838097
$ws.Range("J7").Value = 0.8874158839838097
$ws.Range("M7").Value = 2.350402666666667
$ws.Range("N7").Value = 7.051208000000001
$ws.Range("O7").Value = 0.0613235169202299
$ws.Range("P7").Value = 0.06132351692022989
$ws.Range("Q7").Value = 37.31298000784978
$ws.Range("S7").Value = 0.05441946297676192
$ws.Range("T7").Value = 0.05441946297676192

$ws.Range("I8").Value = 0.07915756308878232
$ws.Range("J8").Value = 0.07915756308878232
$ws.Range("M8").Value = 28.19948866666667
$ws.Range("N8").Value = 84.598466
$ws.Range("O8").Value = 0.7357427920402423
$ws.Range("P8").Value = 0.7357427920402422
$ws.Range("Q8").Value = 39.93228071927467
$ws.Range("R8").Value = 359.3905264734721
$ws.Range("S8").Value = 0.05823960647804233
$ws.Range("T8").Value = 0.05823960647804232

$ws.Range("I9").Value = 0.07915756308878232
$ws.Range("J9").Value = 0.07915756308878232
$ws.Range("O9").Value = 0.2029336910395279
$ws.Range("P9").Value = 0.2029336910395278
$ws.Range("R9").Value = 99.12763923878401
$ws.Range("S9").Value = 0.01606373645130088
$ws.Range("T9").Value = 0.01606373645130088

$ws.Range("I10").Value = 0.07915756308878232
$ws.Range("J10").Value = 0.07915756308878232
$ws.Range("M10").Value = 2.350402666666667
$ws.Range("N10").Value = 7.051208000000001
$ws.Range("O10").Value = 0.0613235169202299
$ws.Range("P10").Value = 0.06132351692022989
$ws.Range("Q10").Value = 3.328320601770668
$ws.Range("R10").Value = 29.954885415936
$ws.Range("S10").Value = 0.004854220159439108
$ws.Range("T10").Value = 0.004854220159439107
